$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update username/password values on row 2 ---
$ws.Range("A2").Value = "ravinskanni@gmail.com"
$ws.Range("B2").Value = "abcd@123E"

# --- Add a duplicate row 3 with the same username/password ---
$ws.Range("A3").Value = "ravinskanni@gmail.com"
$ws.Range("B3").Value = "abcd@123E"

# Give row 3 the same base formatting as row 2 before turning A3 into a
# hyperlink, so the hyperlink style picks up the same vertical-center
# alignment that row 2 already has.
$ws.Range("A2:B2").Copy()
$ws.Range("A3:B3").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- Turn the email cells into mailto hyperlinks ---
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:ravinskanni@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:ravinskanni@gmail.com")

# --- Update the saved selection to match the new active cell ---
[void]$ws.Range("A3").Select()
